$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -5.121723
$ws.Range("D2").Value = 0.000000

# Row 3
$ws.Range("B3").Value = 0.000735
$ws.Range("C3").Value = -5.121723

# Row 4
$ws.Range("B4").Value = 0.005020
$ws.Range("C4").Value = -5.121723
$ws.Range("D4").Value = 0.000000

# Row 5
$ws.Range("B5").Value = 0.009305
$ws.Range("C5").Value = -5.121723

# Row 6
$ws.Range("B6").Value = 0.013590
$ws.Range("C6").Value = -5.170796

# Row 7
$ws.Range("B7").Value = 0.017875
$ws.Range("C7").Value = -5.465231
$ws.Range("D7").Value = 0.343507

# Row 8
$ws.Range("B8").Value = 0.022160
$ws.Range("C8").Value = -5.759665
$ws.Range("D8").Value = 0.637942

# Row 9 removed entirely
$ws.Rows.Item(9).Delete()
